$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Update column C ("Förändrad" / last-changed date) for all data rows 2..232
#    from serial 45192 (2023-09-23) to 45202 (2023-10-03).
for ($r = 2; $r -le 232; $r++) {
    $ws.Cells.Item($r, 3).Value = 45202
}

# 2. Row 232 picks up an explicit custom row height in the saved file.
$ws.Rows.Item(232).RowHeight = 15

# 3. Append the new cleavage-notification row (233).
$ws.Range("A233").Value = "A 46316-2023"

$ws.Range("B233").Value = 45197
$ws.Range("B233").NumberFormat = "YYYY-MM-DD"

$ws.Range("C233").Value = 45202
$ws.Range("C233").NumberFormat = "YYYY-MM-DD"

$ws.Range("D233").Value = "ÖSTERGÖTLANDS LÄN"
$ws.Range("E233").Value = "BOXHOLM"

$ws.Range("G233").Value = 2.3
$ws.Range("H233").Value = 0
$ws.Range("I233").Value = 0
$ws.Range("J233").Value = 0
$ws.Range("K233").Value = 0
$ws.Range("L233").Value = 0
$ws.Range("M233").Value = 0
$ws.Range("N233").Value = 0
$ws.Range("O233").Value = 0
$ws.Range("P233").Value = 0
$ws.Range("Q233").Value = 0

$ws.Range("R233").Value = ""
$ws.Range("R233").WrapText = $true
